# Adds the 2024 row (Decreto 4749/2023) to the "ValoresVenais" table (Tabela4),
# mirroring the author's manual "Add files via upload" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValoresVenais")
$lo = $ws.ListObjects.Item("Tabela4")

# Grow the table by one data row (this also extends the table ref / autofilter
# and the sheet dimension).
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

$values = @(2024, "4749/2023", 326.51, 261.23, 195.88, 156.59, 130.56, 104.51, 91.28, 78.25, 68.5, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 146.23, 324.98, 0, 584.95, 812.42, 1137.43, 0, 243.74, 0, 536.21, 747.42, 0, 0, 219.34, 0, 471.2, 617.43, 0, 0, 0, 3.02)

for ($i = 0; $i -lt $values.Length; $i++) {
    $rng.Cells.Item(1, $i + 1).Value = $values[$i]
}

# Match the author's final selection/scroll position on the sheet.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("X10").Select()
$ws.Range("A44:AN44").Select()
